$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the text values first, in the same order the new shared strings
# were appended in the target workbook (index 19 then index 20).
$ws.Range("C20").Value = "Updating inventory UI showing icon of object."
$ws.Range("C19").Value = "Picking up weapons(if you don’t have it already) and picking up health + showing in inventory. "

# Copy the date/time number formatting from row 18 (format-only paste keeps
# the existing style indices instead of creating new custom number formats).
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("A18").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B18").Copy()
$ws.Range("B20").PasteSpecial(-4122)

# --- Row 19 ---
$ws.Range("A19").Value = Get-Date -Year 2017 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("B19").Value = 0.034722222222222224

# --- Row 20 ---
$ws.Range("A20").Value = Get-Date -Year 2017 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("B20").Value = 0.020833333333333332

# --- Update selection ---
$ws.Range("A21").Select()
